$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vp_sku_list")

# Find the last used row in column A and append the new SKU right below it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = 10041989

# Mirror the activeCell/selection move to the newly added cell, as Excel
# normally does after entering data in the next empty row.
$ws.Cells.Item($newRow, 1).Activate()
